$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105: Especial - new pricing batch (date 44438 -> 44449)
$ws.Cells.Item(105,4).Value = 44449
$ws.Cells.Item(105,13).Value = 100
$ws.Cells.Item(105,14).Value = 10000
$ws.Cells.Item(105,15).Value = 11000
$ws.Cells.Item(105,16).Value = 10500
$ws.Cells.Item(105,19).Value = 656

# Row 106: Primera (date 44438 -> 44449)
$ws.Cells.Item(106,4).Value = 44449
$ws.Cells.Item(106,14).Value = 9500
$ws.Cells.Item(106,15).Value = 9600
$ws.Cells.Item(106,16).Value = 9550
$ws.Cells.Item(106,19).Value = 597

# Row 107: only date changes (44442 -> 44438)
$ws.Cells.Item(107,4).Value = 44438

# Row 108: only date changes (44442 -> 44438)
$ws.Cells.Item(108,4).Value = 44438

# Row 109: date (44435 -> 44442) and quantity M (540 -> 120)
$ws.Cells.Item(109,4).Value = 44442
$ws.Cells.Item(109,13).Value = 120

# Row 110: date (44435 -> 44442) and quantity M (400 -> 60)
$ws.Cells.Item(110,4).Value = 44442
$ws.Cells.Item(110,13).Value = 60

# Row 111: date, category, quantity and price changes
$ws.Cells.Item(111,4).Value = 44435
$ws.Cells.Item(111,12).Value = "Especial"
$ws.Cells.Item(111,13).Value = 540
$ws.Cells.Item(111,14).Value = 9500
$ws.Cells.Item(111,15).Value = 10000
$ws.Cells.Item(111,16).Value = 9750
$ws.Cells.Item(111,19).Value = 609

# Row 112: date, quantity and price changes
$ws.Cells.Item(112,4).Value = 44435
$ws.Cells.Item(112,13).Value = 400
$ws.Cells.Item(112,14).Value = 8500
$ws.Cells.Item(112,16).Value = 8750
$ws.Cells.Item(112,19).Value = 547

# Row 113: date, category, quantity and price changes
$ws.Cells.Item(113,4).Value = 44319
$ws.Cells.Item(113,12).Value = "Primera"
$ws.Cells.Item(113,13).Value = 120
$ws.Cells.Item(113,14).Value = 8500
$ws.Cells.Item(113,15).Value = 9000
$ws.Cells.Item(113,16).Value = 8750
$ws.Cells.Item(113,19).Value = 547

# Row 114: date, quantity and price changes
$ws.Cells.Item(114,4).Value = 44279
$ws.Cells.Item(114,13).Value = 120
$ws.Cells.Item(114,14).Value = 8000
$ws.Cells.Item(114,15).Value = 9000
$ws.Cells.Item(114,16).Value = 8500
$ws.Cells.Item(114,19).Value = 531

# Row 115: date, category, quantity and price changes
$ws.Cells.Item(115,4).Value = 44448
$ws.Cells.Item(115,12).Value = "Especial"
$ws.Cells.Item(115,13).Value = 60
$ws.Cells.Item(115,14).Value = 10000
$ws.Cells.Item(115,15).Value = 11000
$ws.Cells.Item(115,16).Value = 10500
$ws.Cells.Item(115,19).Value = 656

# Row 116: date, category, quantity and price changes
$ws.Cells.Item(116,4).Value = 44448
$ws.Cells.Item(116,12).Value = "Primera"
$ws.Cells.Item(116,13).Value = 60
$ws.Cells.Item(116,14).Value = 9500
$ws.Cells.Item(116,15).Value = 9600
$ws.Cells.Item(116,16).Value = 9550
$ws.Cells.Item(116,19).Value = 597

# Row 117: only date changes (44400 -> 44399)
$ws.Cells.Item(117,4).Value = 44399

# New rows 118 and 119, built from a copy of row 117 formatting, then updated
$ws.Range("A117:T117").Copy($ws.Range("A118:T118"))
$ws.Range("A117:T117").Copy($ws.Range("A119:T119"))

# Row 118: Segunda
$ws.Cells.Item(118,4).Value = 44399
$ws.Cells.Item(118,12).Value = "Segunda"
$ws.Cells.Item(118,13).Value = 120
$ws.Cells.Item(118,14).Value = 8000
$ws.Cells.Item(118,15).Value = 8500
$ws.Cells.Item(118,16).Value = 8250
$ws.Cells.Item(118,19).Value = 516

# Row 119: Primera
$ws.Cells.Item(119,4).Value = 44400
$ws.Cells.Item(119,12).Value = "Primera"
$ws.Cells.Item(119,13).Value = 120
$ws.Cells.Item(119,14).Value = 9500
$ws.Cells.Item(119,15).Value = 10000
$ws.Cells.Item(119,16).Value = 9750
$ws.Cells.Item(119,19).Value = 609
